$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item($row, 1).Value = 42611.883750000001

$ws.Cells.Item($row, 2).Value = 10
$ws.Cells.Item($row, 3).Value = 54
$ws.Cells.Item($row, 4).Value = 43
$ws.Cells.Item($row, 5).Value = 66
$ws.Cells.Item($row, 6).Value = 33
$ws.Cells.Item($row, 7).Value = 19284
$ws.Cells.Item($row, 8).Value = 18995
$ws.Cells.Item($row, 9).Value = 3210
$ws.Cells.Item($row, 10).Value = 388
$ws.Cells.Item($row, 11).Value = 311
$ws.Cells.Item($row, 12).Value = 12
$ws.Cells.Item($row, 13).Value = 6
$ws.Cells.Item($row, 14).Value = "Noun"
